$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.545.53'
$ws.Range('E2').Value = '  -2.57%  '
$ws.Range('D3').Value = '2.305.16'
$ws.Range('E3').Value = '  -3.36%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '540.45'
$ws.Range('E5').Value = '  -1.96%  '
$ws.Range('D6').Value = '127.46'
$ws.Range('E6').Value = '  -5.82%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = '0.568'
$ws.Range('E8').Value = '  -4.13%  '
$ws.Range('D9').Value = '2.302.95'
$ws.Range('E9').Value = '  -3.36%  '
$ws.Range('D10').Value = '0.100'
$ws.Range('E10').Value = '  -1.76%  '
$ws.Range('D11').Value = '5.52'
$ws.Range('E11').Value = '  -0.90%  '
$ws.Range('D12').Value = '0.149'
$ws.Range('E12').Value = '  -1.01%  '
$ws.Range('E13').Value = '  -3.08%  '
$ws.Range('D14').Value = '2.717.51'
$ws.Range('E14').Value = '  -3.30%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').Value = '23.05'
$ws.Range('E15').Value = '  -5.53%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '59.559.86'
$ws.Range('E16').Value = '  -2.38%  '
$ws.Range('E17').Value = '  -3.07%  '
$ws.Range('D18').Value = '2.312.52'
$ws.Range('E18').Value = '  -4.91%  '
$ws.Range('E19').Value = '  -4.55%  '
$ws.Range('E20').Value = '  -5.95%  '
$ws.Range('D21').Value = '309.59'
$ws.Range('E21').Value = '  -3.52%  '
$ws.Range('D22').Value = '6.50'
$ws.Range('E22').Value = '  -6.29%  '
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.90%  '
$ws.Range('D24').Value = '62.98'
$ws.Range('E24').Value = '  -1.22%  '
$ws.Range('E25').Value = '  -3.84%  '
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.35%  '
$ws.Range('D27').Value = '7.69'
$ws.Range('E27').Value = '  -6.81%  '
$ws.Range('D28').Value = '1.33'
$ws.Range('E28').Value = '  -1.93%  '
$ws.Range('E29').Value = '  +3.12%  '
$ws.Range('D30').Value = '171.71'
$ws.Range('E30').Value = '  -0.15%  '
$ws.Range('E31').Value = '  -3.74%  '
$ws.Range('D32').Value = '0.0₃0712'
$ws.Range('E32').Value = '  -6.21%  '
$ws.Range('E33').Value = '  -3.51%  '
$ws.Range('D34').Value = '0.375'
$ws.Range('E34').Value = '  -3.82%  '
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('D36').Value = '1.32'
$ws.Range('E36').Value = '  -8.05%  '
$ws.Range('D37').Value = '17.66'
$ws.Range('E37').Value = '  -2.95%  '
$ws.Range('E38').Value = '  +0.12%  '
$ws.Range('D39').Value = '3.97'
$ws.Range('E39').Value = '  -6.60%  '
$ws.Range('D40').Value = '311.60'
$ws.Range('E40').Value = '  -5.32%  '
$ws.Range('D41').Value = '37.55'
$ws.Range('E41').Value = '  -2.54%  '
$ws.Range('D42').Value = '1.49'
$ws.Range('E42').Value = '  -6.09%  '
$ws.Range('D43').Value = '136.03'
$ws.Range('E43').Value = '  -7.48%  '
$ws.Range('D44').Value = '3.40'
$ws.Range('E44').Value = '  -3.31%  '
$ws.Range('D45').Value = '0.0936'
$ws.Range('E45').Value = '  -2.57%  '
$ws.Range('E46').Value = '  -0.39%  '
$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').Value = '0.0488'
$ws.Range('E47').Value = '  -3.57%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = '18.41'
$ws.Range('E48').Value = '  -6.84%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0223'
$ws.Range('E49').Value = '  +22.49%  '
$ws.Range('E50').Value = '  -2.42%  '
$ws.Range('D51').Value = '10.99'
$ws.Range('E51').Value = '  -0.46%  '
